$d = $word.ActiveDocument

$d.Content.Find.Execute("<id>p090r_a3</id>", $true, $false, $false, $false, $false, $true, 1, $false, "<id>p090r_3</id>", 2)
$d.Content.Find.Execute("<id>p090v_a1</id>", $true, $false, $false, $false, $false, $true, 1, $false, "<id>p090v_1</id>", 2)
$d.Content.Find.Execute("<id>p090v_a2</id>", $true, $false, $false, $false, $false, $true, 1, $false, "<id>p090v_2</id>", 2)
